$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting of the current last row (row 19) down to the new row 20
# before touching row 19's own formatting, so row 20 ends up visually
# identical to what row 19 looked like.
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add the new work log entry as row 20
$ws.Range("A20").Value = 45707
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = "Python file was updated to include more commands, Metasploit framework tool was tested"

# Row 19 is no longer the last row, so it loses its "last row" fill flag
$ws.Range("B19:C19").Interior.Pattern = -4142  # xlPatternNone

# Selection ends up on B21 after data entry
$ws.Range("B21").Select()
